$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Tue Jan 28 22:11:39 EST 2025"
$ws.Range("B3").Value = "Tue Jan 28 22:11:52 EST 2025"
$ws.Range("B4").Value = "Tue Jan 28 22:12:05 EST 2025"
